# New data format: add a 4th column "Erste_Hausnummer" (first/valid house
# number) to the Locations table next to PLZ / Ort / Strasse.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 4).Value = "Erste_Hausnummer"

# Size the new column to fit its header, as Excel does when a user
# auto-fits a freshly typed column.
$ws.Columns("D:D").AutoFit()

# Leave the cursor where the author ended up after populating the new
# column, matching the saved view/selection state.
$ws.Range("D7").Select()
